$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.059.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "'1.649.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").Value = "'217.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("D6").Value = "'0.5209"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.76%  "
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").Value = "'0.2613"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("D9").Value = "'0.06279"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("D10").Value = "'20.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").Value = "'0.07791"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").Value = "'4.473"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("D13").Value = "'1.560.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.44%  "
$ws.Range("D14").Value = "'1.876.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").Value = "'0.5524"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "'0.0₅7981"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.87%  "
$ws.Range("D17").Value = "'64.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("D18").Value = "'26.052.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").Value = "'1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("D20").Value = "'4.620"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.35%  "
$ws.Range("D21").Value = "'193.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  -1.61%  "
$ws.Range("D23").Value = "'5.939"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.68%  "
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").Value = "'146.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").Value = "'0.1198"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.99%  "
$ws.Range("D27").Value = "'7.164"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").Value = "'15.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.53%  "
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").Value = "'0.05590"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.52%  "
$ws.Range("D31").Value = "'1.264"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.56%  "
$ws.Range("D32").Value = "'3.475"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.91%  "
$ws.Range("E33").Value = "  +2.29%  "
$ws.Range("D34").Value = "'1.586"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("D36").Value = "'0.9465"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.76%  "
$ws.Range("D38").Value = "'0.5624"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.21%  "
$ws.Range("D39").Value = "'0.01582"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("D40").Value = "'5.951"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.23%  "
$ws.Range("D41").Value = "'1.057.37"
$ws.Range("D41").Style = "Normal"
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("D43").Value = "'0.8397"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.28%  "
$ws.Range("D44").Value = "'102.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.86%  "
$ws.Range("D45").Value = "'1.788.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("D46").Value = "'56.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.0₈105"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.007"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05345"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.51%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.4331"
$ws.Range("D50").Style = "Normal"
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.903"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.17%  "
